$wb = $excel.ActiveWorkbook

# Values to write into row 11 (columns C:I) for both worksheets.
$values = @(
    0.9625014922111745,
    -0.10000000000000009,
    -0.15982382752014246,
    -0.08000000000000007,
    1.5829618029997903,
    16.12947350163202,
    1.52979216327803
)

foreach ($ws in $wb.Worksheets) {
    $ws.Range("C11").Value = $values[0]
    $ws.Range("D11").Value = $values[1]
    $ws.Range("E11").Value = $values[2]
    $ws.Range("F11").Value = $values[3]
    $ws.Range("G11").Value = $values[4]
    $ws.Range("H11").Value = $values[5]
    $ws.Range("I11").Value = $values[6]
}
